# "High - Low Heat" edit
# Splits the generic "heat" output of the electrolyzer / fuel-synthesizer
# units into "heat_low" / "heat_high", and adds a new "Auxilliary" unit
# (heat_split) on the Units sheet that recombines them into
# "internal_heat". Also updates the Connections sheet so the district
# heat pipeline now draws on "heat_low".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Units" (sheet1 / Table1, range A1:AH8 before the edit)
# ---------------------------------------------------------------
$units = $wb.Worksheets.Item("Units")

# electrolyzer (row 3): Output2 heat -> heat_low
$units.Cells.Item(3, 6).Value = "heat_low"

# fuel_synthesizer (row 7): Output2 heat -> heat_high
$units.Cells.Item(7, 6).Value = "heat_high"

# New row 9: heat_split / Auxilliary unit recombining the heat streams
$units.Cells.Item(9, 1).Value = "heat_split"
$units.Cells.Item(9, 2).Value = "Auxilliary"
$units.Cells.Item(9, 3).Value = "heat_high"
$units.Cells.Item(9, 5).Value = "internal_heat"
$units.Cells.Item(9, 6).Value = "heat_low"
$units.Cells.Item(9, 21).Value = 0.4

# Match the resolution_output / demand column formatting (right aligned,
# same as AF2:AF8 / AG2:AG8) on the new row's equivalent cells (AI9/AJ9)
$units.Cells.Item(9, 35).HorizontalAlignment = -4152
$units.Cells.Item(9, 36).HorizontalAlignment = -4152

# Extend the "h, D, W, M, Q, Y" list validation to the new AI9 cell
$units.Range("AI9").Validation.Add(3, 1, 1, "h, D, W, M, Q, Y")

# Tidy up the view: drop the scrolled-away top-left cell / stale
# selection that used to sit on AH8 and reset back to the sheet origin
$units.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$units.Range("A1").Select()

# ---------------------------------------------------------------
# Sheet "Connections" (sheet2 / Table13)
# ---------------------------------------------------------------
$connections = $wb.Worksheets.Item("Connections")

# pl_dh (row 5): Input1 heat -> heat_low
$connections.Cells.Item(5, 3).Value = "heat_low"
